$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 67.91996
$ws.Range("H2").Value = 203.75988
$ws.Range("I2").Value = 0.02375577759132129
$ws.Range("J2").Value = 0.02375577759132129
$ws.Range("M2").Value = 211.2725676666666
$ws.Range("N2").Value = 633.8177029999999
$ws.Range("O2").Value = 0.6324644927232657
$ws.Range("P2").Value = 0.6324644927232657
$ws.Range("Q2").Value = 14349.62434501729
$ws.Range("R2").Value = 129146.6191051556
$ws.Range("S2").Value = 0.01502468582354174
$ws.Range("T2").Value = 0.01502468582354174
$ws.Range("G3").Value = 67.91996
$ws.Range("H3").Value = 203.75988
$ws.Range("I3").Value = 0.02375577759132129
$ws.Range("J3").Value = 0.02375577759132129
$ws.Range("M3").Value = 59.36675400000001
$ws.Range("N3").Value = 178.100262
$ws.Range("O3").Value = 0.1777200152765546
$ws.Range("P3").Value = 0.1777200152765546
$ws.Range("Q3").Value = 4032.18755700984
$ws.Range("R3").Value = 36289.68801308856
$ws.Range("S3").Value = 0.004221877156436053
$ws.Range("T3").Value = 0.004221877156436053
$ws.Range("G4").Value = 67.91996
$ws.Range("H4").Value = 203.75988
$ws.Range("I4").Value = 0.02375577759132129
$ws.Range("J4").Value = 0.02375577759132129
$ws.Range("M4").Value = 0.4593846666666666
$ws.Range("N4").Value = 1.378154
$ws.Range("O4").Value = 0.001375211620595172
$ws.Range("P4").Value = 0.001375211620595172
$ws.Range("Q4").Value = 31.20138818461333
$ws.Range("R4").Value = 280.81249366152
$ws.Range("S4").Value = 0.00003266922139985943
$ws.Range("T4").Value = 0.00003266922139985943
$ws.Range("G5").Value = 67.91996
$ws.Range("H5").Value = 203.75988
$ws.Range("I5").Value = 0.02375577759132129
$ws.Range("J5").Value = 0.02375577759132129
$ws.Range("M5").Value = 62.94782133333333
$ws.Range("N5").Value = 188.843464
$ws.Range("O5").Value = 0.1884402803795846
$ws.Range("P5").Value = 0.1884402803795846
$ws.Range("Q5").Value = 4275.413507047147
$ws.Range("R5").Value = 38478.72156342432
$ws.Range("S5").Value = 0.004476545389943637
$ws.Range("T5").Value = 0.004476545389943637
$ws.Range("I6").Value = 0.9176057312269553
$ws.Range("J6").Value = 0.9176057312269554
$ws.Range("M6").Value = 211.2725676666666
$ws.Range("N6").Value = 633.8177029999999
$ws.Range("O6").Value = 0.6324644927232657
$ws.Range("P6").Value = 0.6324644927232657
$ws.Range("Q6").Value = 554277.6905249411
$ws.Range("R6").Value = 4988499.21472447
$ws.Range("S6").Value = 0.5803530433204176
$ws.Range("T6").Value = 0.5803530433204176
$ws.Range("I7").Value = 0.9176057312269553
$ws.Range("J7").Value = 0.9176057312269554
$ws.Range("M7").Value = 59.36675400000001
$ws.Range("N7").Value = 178.100262
$ws.Range("O7").Value = 0.1777200152765546
$ws.Range("P7").Value = 0.1777200152765546
$ws.Range("Q7").Value = 155749.8338023653
$ws.Range("R7").Value = 1401748.504221288
$ws.Range("S7").Value = 0.1630769045715086
$ws.Range("T7").Value = 0.1630769045715085
$ws.Range("I8").Value = 0.9176057312269553
$ws.Range("J8").Value = 0.9176057312269554
$ws.Range("M8").Value = 0.4593846666666666
$ws.Range("N8").Value = 1.378154
$ws.Range("O8").Value = 0.001375211620595172
$ws.Range("P8").Value = 0.001375211620595172
$ws.Range("Q8").Value = 1205.204608031767
$ws.Range("R8").Value = 10846.8414722859
$ws.Range("S8").Value = 0.001261902064708039
$ws.Range("T8").Value = 0.001261902064708039
$ws.Range("I9").Value = 0.9176057312269553
$ws.Range("J9").Value = 0.9176057312269554
$ws.Range("M9").Value = 62.94782133333333
$ws.Range("N9").Value = 188.843464
$ws.Range("O9").Value = 0.1884402803795846
$ws.Range("P9").Value = 0.1884402803795846
$ws.Range("Q9").Value = 165144.8336031249
$ws.Range("R9").Value = 1486303.502428124
$ws.Range("S9").Value = 0.1729138812703212
$ws.Range("T9").Value = 0.1729138812703212
$ws.Range("G10").Value = 1.376679
$ws.Range("H10").Value = 4.130037
$ws.Range("I10").Value = 0.0004815091195378001
$ws.Range("J10").Value = 0.0004815091195378002
$ws.Range("M10").Value = 211.2725676666666
$ws.Range("N10").Value = 633.8177029999999
$ws.Range("O10").Value = 0.6324644927232657
$ws.Range("P10").Value = 0.6324644927232657
$ws.Range("Q10").Value = 290.8545071827789
$ws.Range("R10").Value = 2617.69056464501
$ws.Range("S10").Value = 0.000304537421030101
$ws.Range("T10").Value = 0.000304537421030101
$ws.Range("G11").Value = 1.376679
$ws.Range("H11").Value = 4.130037
$ws.Range("I11").Value = 0.0004815091195378001
$ws.Range("J11").Value = 0.0004815091195378002
$ws.Range("M11").Value = 59.36675400000001
$ws.Range("N11").Value = 178.100262
$ws.Range("O11").Value = 0.1777200152765546
$ws.Range("P11").Value = 0.1777200152765546
$ws.Range("Q11").Value = 81.72896352996601
$ws.Range("R11").Value = 735.5606717696941
$ws.Range("S11").Value = 0.00008557380808005818
$ws.Range("T11").Value = 0.00008557380808005818
$ws.Range("G12").Value = 1.376679
$ws.Range("H12").Value = 4.130037
$ws.Range("I12").Value = 0.0004815091195378001
$ws.Range("J12").Value = 0.0004815091195378002
$ws.Range("M12").Value = 0.4593846666666666
$ws.Range("N12").Value = 1.378154
$ws.Range("O12").Value = 0.001375211620595172
$ws.Range("P12").Value = 0.001375211620595172
$ws.Range("Q12").Value = 0.6324252235219999
$ws.Range("R12").Value = 5.691827011697999
$ws.Range("S12").Value = 0.0000006621769366109325
$ws.Range("T12").Value = 0.0000006621769366109326
$ws.Range("G13").Value = 1.376679
$ws.Range("H13").Value = 4.130037
$ws.Range("I13").Value = 0.0004815091195378001
$ws.Range("J13").Value = 0.0004815091195378002
$ws.Range("M13").Value = 62.94782133333333
$ws.Range("N13").Value = 188.843464
$ws.Range("O13").Value = 0.1884402803795846
$ws.Range("P13").Value = 0.1884402803795846
$ws.Range("Q13").Value = 86.65894372535199
$ws.Range("R13").Value = 779.9304935281679
$ws.Range("S13").Value = 0.00009073571349102997
$ws.Range("T13").Value = 0.00009073571349102997
$ws.Range("G14").Value = 164.8447596666666
$ws.Range("H14").Value = 494.534279
$ws.Range("I14").Value = 0.05765632735555414
$ws.Range("J14").Value = 0.05765632735555416
$ws.Range("M14").Value = 211.2725676666666
$ws.Range("N14").Value = 633.8177029999999
$ws.Range("O14").Value = 0.6324644927232657
$ws.Range("P14").Value = 0.6324644927232657
$ws.Range("Q14").Value = 34827.17564117123
$ws.Range("R14").Value = 313444.5807705411
$ws.Range("S14").Value = 0.0364655798332171
$ws.Range("T14").Value = 0.0364655798332171
$ws.Range("G15").Value = 164.8447596666666
$ws.Range("H15").Value = 494.534279
$ws.Range("I15").Value = 0.05765632735555414
$ws.Range("J15").Value = 0.05765632735555416
$ws.Range("M15").Value = 59.36675400000001
$ws.Range("N15").Value = 178.100262
$ws.Range("O15").Value = 0.1777200152765546
$ws.Range("P15").Value = 0.1777200152765546
$ws.Range("Q15").Value = 9786.298295320123
$ws.Range("R15").Value = 88076.68465788109
$ws.Range("S15").Value = 0.01024668337841911
$ws.Range("T15").Value = 0.01024668337841911
$ws.Range("G16").Value = 164.8447596666666
$ws.Range("H16").Value = 494.534279
$ws.Range("I16").Value = 0.05765632735555414
$ws.Range("J16").Value = 0.05765632735555416
$ws.Range("M16").Value = 0.4593846666666666
$ws.Range("N16").Value = 1.378154
$ws.Range("O16").Value = 0.001375211620595172
$ws.Range("P16").Value = 0.001375211620595172
$ws.Range("Q16").Value = 75.72715497121843
$ws.Range("R16").Value = 681.5443947409659
$ws.Range("S16").Value = 0.00007928965138019737
$ws.Range("T16").Value = 0.00007928965138019739
$ws.Range("G17").Value = 164.8447596666666
$ws.Range("H17").Value = 494.534279
$ws.Range("I17").Value = 0.05765632735555414
$ws.Range("J17").Value = 0.05765632735555416
$ws.Range("M17").Value = 62.94782133333333
$ws.Range("N17").Value = 188.843464
$ws.Range("O17").Value = 0.1884402803795846
$ws.Range("P17").Value = 0.1884402803795846
$ws.Range("Q17").Value = 10376.6184792336
$ws.Range("R17").Value = 93389.56631310245
$ws.Range("S17").Value = 0.01086477449253774
$ws.Range("T17").Value = 0.01086477449253774
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1.431418
$ws.Range("H18").Value = 4.294254
$ws.Range("I18").Value = 0.0005006547066313635
$ws.Range("J18").Value = 0.0005006547066313636
$ws.Range("M18").Value = 211.2725676666666
$ws.Range("N18").Value = 633.8177029999999
$ws.Range("O18").Value = 0.6324644927232657
$ws.Range("P18").Value = 0.6324644927232657
$ws.Range("Q18").Value = 302.4193562642847
$ws.Range("R18").Value = 2721.774206378562
$ws.Range("S18").Value = 0.0003166463250591207
$ws.Range("T18").Value = 0.0003166463250591208
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1.431418
$ws.Range("H19").Value = 4.294254
$ws.Range("I19").Value = 0.0005006547066313635
$ws.Range("J19").Value = 0.0005006547066313636
$ws.Range("M19").Value = 59.36675400000001
$ws.Range("N19").Value = 178.100262
$ws.Range("O19").Value = 0.1777200152765546
$ws.Range("P19").Value = 0.1777200152765546
$ws.Range("Q19").Value = 84.97864027717202
$ws.Range("R19").Value = 764.8077624945481
$ws.Range("S19").Value = 0.00008897636211080487
$ws.Range("T19").Value = 0.00008897636211080488
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1.431418
$ws.Range("H20").Value = 4.294254
$ws.Range("I20").Value = 0.0005006547066313635
$ws.Range("J20").Value = 0.0005006547066313636
$ws.Range("M20").Value = 0.4593846666666666
$ws.Range("N20").Value = 1.378154
$ws.Range("O20").Value = 0.001375211620595172
$ws.Range("P20").Value = 0.001375211620595172
$ws.Range("Q20").Value = 0.6575714807906666
$ws.Range("R20").Value = 5.918143327116
$ws.Range("S20").Value = 0.0000006885061704651179
$ws.Range("T20").Value = 0.000000688506170465118
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1.431418
$ws.Range("H21").Value = 4.294254
$ws.Range("I21").Value = 0.0005006547066313635
$ws.Range("J21").Value = 0.0005006547066313636
$ws.Range("M21").Value = 62.94782133333333
$ws.Range("N21").Value = 188.843464
$ws.Range("O21").Value = 0.1884402803795846
$ws.Range("P21").Value = 0.1884402803795846
$ws.Range("Q21").Value = 90.10464451731734
$ws.Range("R21").Value = 810.941800655856
$ws.Range("S21").Value = 0.00009434351329097282
$ws.Range("T21").Value = 0.00009434351329097282
